$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (strikeouts) values for column G, rows 2-17,
# regenerated from the source stats instead of the old placeholder "Strike#" data.
$newK = @{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 0
    7  = 1
    8  = 4
    9  = 1
    10 = 0
    11 = 1
    12 = 1
    13 = 2
    14 = 2
    15 = 2
    16 = 2
    17 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
